$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.339.03"
$ws.Range("E2").Value = "  -4.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.340.76"
$ws.Range("E3").Value = "  -2.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.59"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.73"
$ws.Range("E6").Value = "  +3.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.340.08"
$ws.Range("E8").Value = "  -2.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  -0.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.120"
$ws.Range("E11").Value = "  -2.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  +0.28%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.911.67"
$ws.Range("E13").Value = "  -2.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.120"
$ws.Range("E14").Value = "  +0.44%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  -2.06%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.345.14"
$ws.Range("E16").Value = "  -2.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.80"
$ws.Range("E17").Value = "  -0.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.489.81"
$ws.Range("E18").Value = "  -4.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.70"
$ws.Range("E19").Value = "  +3.74%  "

# Row 20
$ws.Range("E20").Value = "  +1.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.26"
$ws.Range("E21").Value = "  -5.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.39"
$ws.Range("E22").Value = "  -2.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").Value = "  -0.26%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.472.12"
$ws.Range("E25").Value = "  -2.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.68"
$ws.Range("E26").Value = "  -4.43%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +4.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.67"
$ws.Range("E28").Value = "  +18.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").Value = "  +7.41%  "

# Row 30
$ws.Range("E30").Value = "  -0.30%  "

# Row 31
$ws.Range("E31").Value = "  +1.98%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -1.74%  "

# Row 33
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  +0.50%  "

# Row 34
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.372.15"
$ws.Range("E35").Value = "  -2.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.99"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.48"
$ws.Range("E37").Value = "  +2.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.94"
$ws.Range("E38").Value = "  +2.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +0.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.96"
$ws.Range("E40").Value = "  -0.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0778"
$ws.Range("E41").Value = "  +1.51%  "

# Row 42
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.21"
$ws.Range("E43").Value = "  +10.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").Value = "  +2.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.22"
$ws.Range("E45").Value = "  -2.23%  "

# Row 46
$ws.Range("E46").Value = "  -3.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.75"
$ws.Range("E47").Value = "  +3.31%  "

# Row 48
$ws.Range("E48").Value = "  -0.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.91"
$ws.Range("E49").Value = "  +2.95%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.53"
$ws.Range("E50").Value = "  +10.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.897"
$ws.Range("E51").Value = "  +3.00%  "
